$wb = $excel.ActiveWorkbook

# --- Sheet: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("D15").Value = 1299.46
$ws1.Range("M15").Value = 1196.16
$ws1.Range("E19").Value = 145.95
$ws1.Range("M19").Value = 349.03
$ws1.Range("E34").Value = "1 de 32"
$ws1.Range("M34").Value = "3 de 32"

# --- Sheet: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("G2").Value = 3000
$ws2.Range("G4").Value = 1500
$ws2.Range("G8").Value = 0
$ws2.Range("G9").Value = 0
$ws2.Range("G10").Value = 500
$ws2.Range("G12").Value = 2000
$ws2.Range("F15").Value = 2495.62
$ws2.Range("G15").Value = 4000
$ws2.Range("G16").Value = 500
$ws2.Range("G17").Value = 1000
$ws2.Range("G18").Value = 0
$ws2.Range("F19").Value = 494.98
$ws2.Range("G19").Value = 4000
$ws2.Range("G26").Value = 900
$ws2.Range("G27").Value = 8000
$ws2.Range("G29").Value = 1500
$ws2.Range("G31").Value = 1500
$ws2.Range("G33").Value = 0
$ws2.Range("F34").Value = 3088.42
$ws2.Range("G34").Value = 33900
